$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.769.30"
$ws.Range("E2").Value = "  +1.84%  "
$ws.Range("D3").Value = "3.560.58"
$ws.Range("E3").Value = "  +1.59%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'585.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").Value = "'188.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.57%  "
$ws.Range("D7").Value = "3.554.47"
$ws.Range("E7").Value = "  +1.68%  "
$ws.Range("D8").Value = "'0.622"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.70%  "
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("E10").Value = "  +6.40%  "
$ws.Range("D11").Value = "'0.644"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").Value = "'53.95"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").Value = "'0.0000308"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.76%  "
$ws.Range("D14").Value = "'9.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("D15").Value = "4.126.70"
$ws.Range("E15").Value = "  +1.52%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.615.30"
$ws.Range("E16").Value = "  +3.43%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "70.780.13"
$ws.Range("E17").Value = "  +1.91%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'12.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.79%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'18.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.81%  "
$ws.Range("D20").Value = "'566.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.06%  "
$ws.Range("E21").Value = "  +0.72%  "
$ws.Range("D22").Value = "'0.993"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.74%  "
$ws.Range("E23").Value = "  -3.94%  "
$ws.Range("D24").Value = "'4.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.08%  "
$ws.Range("D25").Value = "'4.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").Value = "'93.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.29%  "
$ws.Range("D27").Value = "'11.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("E28").Value = "  -2.48%  "
$ws.Range("D29").Value = "'9.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.88%  "
$ws.Range("E30").Value = "  +1.50%  "
$ws.Range("D31").Value = "'7.04"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.13%  "
$ws.Range("D32").Value = "'12.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.81%  "
$ws.Range("D33").Value = "'3.93"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +25.96%  "
$ws.Range("E34").Value = "  +1.37%  "
$ws.Range("D35").Value = "'63.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.32%  "
$ws.Range("E36").Value = "  +5.38%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").Value = "'528.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.75%  "
$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").Value = "'0.406"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.65%  "
$ws.Range("D39").Value = "'38.14"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "3.640.12"
$ws.Range("E40").Value = "  +8.87%  "
$ws.Range("B41").Value = "Dai"
$ws.Range("C41").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D41").Value = "'0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").Value = "0.0₃0785"
$ws.Range("E42").Value = "  +2.92%  "
$ws.Range("E43").Value = "  +4.45%  "
$ws.Range("E44").Value = "  +2.96%  "
$ws.Range("D45").Value = "'0.0457"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.51%  "
$ws.Range("E46").Value = "  -0.94%  "
$ws.Range("E47").Value = "  -1.86%  "
$ws.Range("E48").Value = "  +2.54%  "
$ws.Range("D49").Value = "'9.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.07%  "
$ws.Range("D50").Value = "'0.999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("E51").Value = "  +8.64%  "
